# Auto-generated edit script applying cryptos.xlsx scraped-data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.NumberFormat = "General"
}

Set-TextCell 'D2' '328.82'
Set-TextCell 'E2' '0.53%'
Set-TextCell 'D3' '44.19'
Set-TextCell 'E3' '1.01%'
Set-TextCell 'D4' '5.478'
Set-TextCell 'E4' '-1.32%'
Set-TextCell 'D5' '0.08062'
Set-TextCell 'E5' '0.53%'
Set-TextCell 'D6' '2.044'
Set-TextCell 'E6' '7.61%'
Set-TextCell 'D7' '0.9535'
Set-TextCell 'E7' '0.80%'
Set-TextCell 'D8' '0.1128'
Set-TextCell 'E8' '-3.92%'
Set-TextCell 'D9' '0.1874'
Set-TextCell 'E9' '1.94%'
Set-TextCell 'D10' '10.26'
Set-TextCell 'E10' '-2.36%'
Set-TextCell 'D11' '0.09919'
Set-TextCell 'E11' '3.22%'
Set-TextCell 'D12' '0.04778'
Set-TextCell 'E12' '6.93%'
Set-TextCell 'D13' '0.1062'
Set-TextCell 'E13' '-0.33%'
Set-TextCell 'D14' '0.001270'
Set-TextCell 'E14' '-1.34%'
Set-TextCell 'D15' '0.04084'
Set-TextCell 'E15' '-2.67%'
Set-TextCell 'D16' '0.005844'
Set-TextCell 'E16' '-2.47%'
Set-TextCell 'B17' 'LEO'
Set-TextCell 'C17' 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextCell 'D17' '3.376'
Set-TextCell 'E17' '-0.88%'
Set-TextCell 'B18' 'GateToken'
Set-TextCell 'C18' 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextCell 'D18' '4.416'
Set-TextCell 'E18' '3.54%'
Set-TextCell 'B19' 'BTSEToken'
Set-TextCell 'C19' 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextCell 'D19' '2.622'
Set-TextCell 'E19' '3.24%'
Set-TextCell 'B20' 'BitpandaEcosystemToken'
Set-TextCell 'C20' 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
Set-TextCell 'D20' '0.3409'
Set-TextCell 'E20' '-1.04%'
Set-TextCell 'B21' 'ProBitToken'
Set-TextCell 'C21' 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
Set-TextCell 'D21' '0.1400'
Set-TextCell 'E21' '0.24%'
Set-TextCell 'B22' 'ZBToken'
Set-TextCell 'C22' 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
Set-TextCell 'D22' '0.2581'
Set-TextCell 'E22' '2.95%'
Set-TextCell 'B23' 'BitKan'
Set-TextCell 'C23' 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
Set-TextCell 'D23' '0.001307'
Set-TextCell 'E23' '4.74%'
Set-TextCell 'B24' 'HotbitToken'
Set-TextCell 'C24' 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
Set-TextCell 'D24' '0.004345'
Set-TextCell 'E24' '1.66%'
Set-TextCell 'E25' '-0.89%'
Set-TextCell 'D26' '0.0003743'
Set-TextCell 'E26' '-6.32%'
Set-TextCell 'D38' '0.02576'
Set-TextCell 'E38' '-2.51%'
Set-TextCell 'D39' '0.05674'
Set-TextCell 'E39' '3.13%'
Set-TextCell 'D40' '0.007716'
Set-TextCell 'E40' '1.85%'
Set-TextCell 'D41' '0.1399'
Set-TextCell 'E41' '0.63%'
Set-TextCell 'D42' '0.007354'
Set-TextCell 'E42' '-9.96%'
Set-TextCell 'D43' '0.002010'
Set-TextCell 'E43' '0.26%'
Set-TextCell 'D44' '0.008526'
Set-TextCell 'E44' '-3.18%'
Set-TextCell 'E45' '1.90%'
Set-TextCell 'E46' '-0.12%'
Set-TextCell 'D47' '0.0005805'
Set-TextCell 'E47' '-0.11%'
Set-TextCell 'D48' '0.003511'
Set-TextCell 'E48' '0.13%'
Set-TextCell 'D49' '0.003501'
Set-TextCell 'E49' '53.98%'
Set-TextCell 'E50' '-0.12%'
Set-TextCell 'E51' '-0.12%'
